# "Change in chrome driver" - update the Selenium/Jira tracking values in
# the OrangeHRM_Excel workbook: the Admin sheet's hashed field id (D2) and
# the Jira sheet's bug-tracking row (B2/A3/B3) move on to their next
# values/ticket.

$wb = $excel.ActiveWorkbook

# --- Admin sheet: D2 hash id moves on to a new value -----------------
$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Range("D2").Value = "1342182A"

# --- Jira sheet: tracking row gets a new ticket/component/id ---------
$wsJira = $wb.Worksheets.Item("Jira")
$wsJira.Range("B2").Value = "C740395689588328E5DA6BCCD7E88F35"
$wsJira.Range("A3").Value = "PersonalDetails"

# B3 ("10600") looks like a plain integer; a direct Value assignment
# would be auto-coerced to a Number by Excel's type inference, but the
# source data keeps it as literal Text (matching the existing "10400").
# Round-trip the literal string through a text formula + values-only
# paste so the stored type stays Text without touching the cell style.
$helper = $wsJira.Range("D100")
$helper.Formula = "=""10600"""
$helper.Copy()
$wsJira.Range("B3").PasteSpecial(-4163)
$helper.Clear()

$excel.CutCopyMode = 0
